$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5152.273
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5152.273
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15456.819
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -17204.819

$ws.Range("H72").Value = 5152.273
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5152.273
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 46370.457
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -55106.457

$ws.Range("H86").Value = 3338
$ws.Range("I86").Value = 3201.2
$ws.Range("J86").Value = 3680
$ws.Range("K86").Value = 3201.2
$ws.Range("L86").Value = 3680
$ws.Range("M86").Value = -2078.2
$ws.Range("N86").Value = -5926

$ws.Range("H89").Value = 3338
$ws.Range("I89").Value = 3201.2
$ws.Range("J89").Value = 3680
$ws.Range("K89").Value = 16006
$ws.Range("L89").Value = 18400
$ws.Range("M89").Value = -10390
$ws.Range("N89").Value = -29632

$ws.Range("H106").Value = 66669236
$ws.Range("I106").Value = 23812124
$ws.Range("K106").Value = 23812124
$ws.Range("M106").Value = -23811493

$ws.Range("H137").Value = 1346.8889
$ws.Range("I137").Value = 1048.2
$ws.Range("J137").Value = 2392.3
$ws.Range("K137").Value = 3144.6
$ws.Range("L137").Value = 7176.900000000001
$ws.Range("M137").Value = -594.6000000000004
$ws.Range("N137").Value = -12276.9

$ws.Range("H138").Value = 2914.1384
$ws.Range("I138").Value = 1306.4062
$ws.Range("J138").Value = 4473.1514
$ws.Range("K138").Value = 3919.2186
$ws.Range("L138").Value = 13419.4542
$ws.Range("M138").Value = 1220.7814
$ws.Range("N138").Value = -23699.4542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1413.8125
$ws.Range("I97").Value = 1308.25
$ws.Range("J97").Value = 1730.5
$ws.Range("K97").Value = 1308.25
$ws.Range("L97").Value = 1730.5
$ws.Range("M97").Value = -812.25
$ws.Range("N97").Value = -2722.5

$ws.Range("H132").Value = 1854177.4
$ws.Range("I132").Value = 1649.6818
$ws.Range("K132").Value = 4949.0454
$ws.Range("M132").Value = -2419.0454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1691.1111
$ws.Range("I86").Value = 1652.7273
$ws.Range("J86").Value = 1751.4286
$ws.Range("K86").Value = 1652.7273
$ws.Range("L86").Value = 1751.4286
$ws.Range("M86").Value = -529.7273
$ws.Range("N86").Value = -3997.4286

$ws.Range("H89").Value = 1691.1111
$ws.Range("I89").Value = 1652.7273
$ws.Range("J89").Value = 1751.4286
$ws.Range("K89").Value = 8263.636500000001
$ws.Range("L89").Value = 8757.143
$ws.Range("M89").Value = -2647.636500000001
$ws.Range("N89").Value = -19989.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7933187
$ws.Range("I31").Value = 1777.871
$ws.Range("J31").Value = 15164766
$ws.Range("K31").Value = 1777.871
$ws.Range("L31").Value = 15164766
$ws.Range("M31").Value = -1482.871
$ws.Range("N31").Value = -15165356

$ws.Range("H34").Value = 7933187
$ws.Range("I34").Value = 1777.871
$ws.Range("J34").Value = 15164766
$ws.Range("K34").Value = 1777.871
$ws.Range("L34").Value = 15164766
$ws.Range("M34").Value = -1575.871
$ws.Range("N34").Value = -15165170

$ws.Range("H58").Value = 2565007.5
$ws.Range("I58").Value = 3205652
$ws.Range("J58").Value = 2429.6924
$ws.Range("K58").Value = 3205652
$ws.Range("L58").Value = 2429.6924
$ws.Range("M58").Value = -3205449
$ws.Range("N58").Value = -2835.6924

$ws.Range("H70").Value = 45000
$ws.Range("I70").Value = 40000
$ws.Range("J70").Value = 50000
$ws.Range("K70").Value = 40000
$ws.Range("L70").Value = 50000
$ws.Range("M70").Value = -39685
$ws.Range("N70").Value = -50630

$ws.Range("H73").Value = 45000
$ws.Range("I73").Value = 40000
$ws.Range("J73").Value = 50000
$ws.Range("K73").Value = 40000
$ws.Range("L73").Value = 50000
$ws.Range("M73").Value = -38908
$ws.Range("N73").Value = -52184

$ws.Range("H107").Value = 12821155
$ws.Range("I107").Value = 18518930
$ws.Range("J107").Value = 1158.875
$ws.Range("K107").Value = 18518930
$ws.Range("L107").Value = 1158.875
$ws.Range("M107").Value = -18517010
$ws.Range("N107").Value = -4998.875

$ws.Range("H132").Value = 2779160
$ws.Range("I132").Value = 3704916.8
$ws.Range("K132").Value = 11114750.4
$ws.Range("M132").Value = -11112220.4

$ws.Range("H134").Value = 5748959
$ws.Range("I134").Value = 10103404
$ws.Range("J134").Value = 1092.44
$ws.Range("K134").Value = 30310212
$ws.Range("L134").Value = 3277.32
$ws.Range("M134").Value = -30307677
$ws.Range("N134").Value = -8347.32

$ws.Range("H136").Value = 2565007.5
$ws.Range("I136").Value = 3205652
$ws.Range("J136").Value = 2429.6924
$ws.Range("K136").Value = 9616956
$ws.Range("L136").Value = 7289.0772
$ws.Range("M136").Value = -9614406
$ws.Range("N136").Value = -12389.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 417249.78
$ws.Range("I113").Value = 664.5454999999999
$ws.Range("J113").Value = 769745
$ws.Range("K113").Value = 1993.6365
$ws.Range("L113").Value = 2309235
$ws.Range("M113").Value = 176.3635000000002
$ws.Range("N113").Value = -2313575

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 39833.332
$ws.Range("J74").Value = 39833.332
$ws.Range("L74").Value = 39833.332
$ws.Range("N74").Value = -41705.332

$ws.Range("H77").Value = 39833.332
$ws.Range("J77").Value = 39833.332
$ws.Range("L77").Value = 119499.996
$ws.Range("N77").Value = -128859.996

$ws.Range("H80").Value = 2783.111
$ws.Range("I80").Value = 2280
$ws.Range("J80").Value = 3412
$ws.Range("K80").Value = 2280
$ws.Range("L80").Value = 3412
$ws.Range("M80").Value = -1282
$ws.Range("N80").Value = -5408

$ws.Range("H83").Value = 2783.111
$ws.Range("I83").Value = 2280
$ws.Range("J83").Value = 3412
$ws.Range("K83").Value = 11400
$ws.Range("L83").Value = 17060
$ws.Range("M83").Value = -6408
$ws.Range("N83").Value = -27044

$ws.Range("H97").Value = 840
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 800
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -304
$ws.Range("N97").Value = -1992

$ws.Range("H132").Value = 5749244.5
$ws.Range("I132").Value = 6174899.5
$ws.Range("K132").Value = 18524698.5
$ws.Range("M132").Value = -18522168.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8152037
$ws.Range("I122").Value = 17875642
$ws.Range("J122").Value = 1669633.4
$ws.Range("K122").Value = 53626926
$ws.Range("L122").Value = 5008900.199999999
$ws.Range("M122").Value = -53624476
$ws.Range("N122").Value = -5013800.199999999

$ws.Range("H132").Value = 13362782
$ws.Range("I132").Value = 20556930
$ws.Range("J132").Value = 2222.7856
$ws.Range("K132").Value = 61670790
$ws.Range("L132").Value = 6668.3568
$ws.Range("M132").Value = -61668260
$ws.Range("N132").Value = -11728.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3086
$ws.Range("I62").Value = 3086
$ws.Range("K62").Value = 3086
$ws.Range("M62").Value = -2462

$ws.Range("H65").Value = 3086
$ws.Range("I65").Value = 3086
$ws.Range("K65").Value = 15430
$ws.Range("M65").Value = -12310

$ws.Range("H100").Value = 370
$ws.Range("I100").Value = 337.5
$ws.Range("J100").Value = 500
$ws.Range("K100").Value = 675
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -134
$ws.Range("N100").Value = -2082

Write-Output "Applied all Leve price/profit updates"
